$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 24.65
$ws.Range("F2").Value = 21.1
$ws.Range("G2").Value = 28.2
$ws.Range("H2").Value = 5.02045814642449
$ws.Range("I2").Value = 0.203669701680507
$ws.Range("E3").Value = 17.52
$ws.Range("F3").Value = 16.5
$ws.Range("G3").Value = 18.5
$ws.Range("H3").Value = 0.920326029187483
$ws.Range("I3").Value = 0.0525300244970024
$ws.Range("E4").Value = 12.8
$ws.Range("F4").Value = 8.7
$ws.Range("G4").Value = 17.1
$ws.Range("H4").Value = 2.44994331000173
$ws.Range("I4").Value = 0.191401821093885
$ws.Range("E5").Value = 34.02
$ws.Range("F5").Value = 28.5
$ws.Range("G5").Value = 39.1
$ws.Range("H5").Value = 3.33826302139301
$ws.Range("I5").Value = 0.0981264850497652
$ws.Range("E6").Value = 24.12
$ws.Range("F6").Value = 19.7
$ws.Range("G6").Value = 28.8
$ws.Range("H6").Value = 3.33992681223872
$ws.Range("I6").Value = 0.138471260872252
$ws.Range("E7").Value = 19.95
$ws.Range("F7").Value = 14.9
$ws.Range("G7").Value = 24.1
$ws.Range("H7").Value = 3.47315098689616
$ws.Range("I7").Value = 0.174092781298053
$ws.Range("E8").Value = 33.13
$ws.Range("F8").Value = 29.4
$ws.Range("G8").Value = 35.5
$ws.Range("H8").Value = 2.11189751434844
$ws.Range("I8").Value = 0.0637457746558538
$ws.Range("E9").Value = 23.0555555555556
$ws.Range("F9").Value = 17.3
$ws.Range("G9").Value = 29.9
$ws.Range("H9").Value = 3.81153220867642
$ws.Range("I9").Value = 0.165319469291989
$ws.Range("E10").Value = 12.2
$ws.Range("F10").Value = 12.2
$ws.Range("G10").Value = 12.2
$ws.Range("E11").Value = 20.76
$ws.Range("F11").Value = 20.1
$ws.Range("G11").Value = 22.2
$ws.Range("H11").Value = 0.844393273303381
$ws.Range("I11").Value = 0.040674049773766
$ws.Range("E12").Value = 22.54
$ws.Range("F12").Value = 18.7
$ws.Range("G12").Value = 30.3
$ws.Range("H12").Value = 3.65489017922868
$ws.Range("I12").Value = 0.16215129455318
$ws.Range("E13").Value = 7.36666666666667
$ws.Range("F13").Value = 6.8
$ws.Range("G13").Value = 7.8
$ws.Range("H13").Value = 0.513160143944688
$ws.Range("I13").Value = 0.0696597480467903
$ws.Range("E14").Value = 9.41
$ws.Range("F14").Value = 8.5
$ws.Range("G14").Value = 10.5
$ws.Range("H14").Value = 0.570477382938574
$ws.Range("I14").Value = 0.0606245890476699
$ws.Range("E15").Value = 9.96666666666667
$ws.Range("F15").Value = 6.8
$ws.Range("G15").Value = 12.9
$ws.Range("H15").Value = 1.58034806292791
$ws.Range("I15").Value = 0.158563350795443
$ws.Range("E16").Value = 20.33
$ws.Range("F16").Value = 13.5
$ws.Range("G16").Value = 22.4
$ws.Range("H16").Value = 2.5060149862104
$ws.Range("I16").Value = 0.123266846345814
$ws.Range("E17").Value = 10.8583333333333
$ws.Range("F17").Value = 9.9
$ws.Range("G17").Value = 12.2
$ws.Range("H17").Value = 0.657071095272364
$ws.Range("I17").Value = 0.0605130709383604
$ws.Range("E18").Value = 8.12
$ws.Range("F18").Value = 7.2
$ws.Range("G18").Value = 8.8
$ws.Range("H18").Value = 0.489444129146071
$ws.Range("I18").Value = 0.0602763705844915
$ws.Range("E19").Value = 12.6875
$ws.Range("G19").Value = 13.5
$ws.Range("H19").Value = 0.458062690656452
$ws.Range("I19").Value = 0.0361034633029716
$ws.Range("E20").Value = 13.66
$ws.Range("F20").Value = 12.5
$ws.Range("G20").Value = 15.3
$ws.Range("H20").Value = 0.732120208708925
$ws.Range("I20").Value = 0.0535959157180765
$ws.Range("E21").Value = 5.5
$ws.Range("F21").Value = 5.5
$ws.Range("G21").Value = 5.5
$ws.Range("E22").Value = 5.1
$ws.Range("F22").Value = 4.5
$ws.Range("G22").Value = 5.7
$ws.Range("H22").Value = 0.496655480858378
$ws.Range("I22").Value = 0.0973834276192898
$ws.Range("E23").Value = 5.3
$ws.Range("F23").Value = 4.8
$ws.Range("G23").Value = 5.9
$ws.Range("H23").Value = 0.556776436283002
$ws.Range("I23").Value = 0.105052157789246
$ws.Range("E24").Value = 6.33333333333333
$ws.Range("G24").Value = 6.7
$ws.Range("H24").Value = 0.294392028877595
$ws.Range("I24").Value = 0.0464829519280413
$ws.Range("E25").Value = 5.125
$ws.Range("F25").Value = 4.9
$ws.Range("G25").Value = 5.4
$ws.Range("H25").Value = 0.206155281280883
$ws.Range("I25").Value = 0.0402254207377333
$ws.Range("E26").Value = 23.9
$ws.Range("F26").Value = 23.9
$ws.Range("G26").Value = 23.9
$ws.Range("E27").Value = 26.25
$ws.Range("F27").Value = 26.1
$ws.Range("G27").Value = 26.4
$ws.Range("H27").Value = 0.212132034355962
$ws.Range("I27").Value = 0.00808122035641761
$ws.Range("E28").Value = 18.3
$ws.Range("G28").Value = 21.4
$ws.Range("H28").Value = 1.9131126469709
$ws.Range("I28").Value = 0.104541674697863
$ws.Range("E29").Value = 14.44
$ws.Range("F29").Value = 9.7
$ws.Range("G29").Value = 20.1
$ws.Range("H29").Value = 3.26333163091546
$ws.Range("I29").Value = 0.225992495215752
$ws.Range("E30").Value = 8.91428571428571
$ws.Range("F30").Value = 5.6
$ws.Range("H30").Value = 3.12699460154559
$ws.Range("I30").Value = 0.350784650814409
$ws.Range("E31").Value = 23.5333333333333
$ws.Range("F31").Value = 22.5
$ws.Range("G31").Value = 24.2
$ws.Range("H31").Value = 0.907377172587746
$ws.Range("I31").Value = 0.0385571036510374
$ws.Range("E32").Value = 7.74444444444444
$ws.Range("F32").Value = 6.9
$ws.Range("H32").Value = 0.705533682950558
$ws.Range("I32").Value = 0.0911019102805598
$ws.Range("E33").Value = 40.8
$ws.Range("F33").Value = 40.8
$ws.Range("G33").Value = 40.8
$ws.Range("E34").Value = 20.55
$ws.Range("F34").Value = 15.3
$ws.Range("G34").Value = 23.6
$ws.Range("H34").Value = 2.6336497700171
$ws.Range("I34").Value = 0.128158139660199
$ws.Range("E35").Value = 22.11
$ws.Range("F35").Value = 18.8
$ws.Range("G35").Value = 26.2
$ws.Range("H35").Value = 2.84310081112546
$ws.Range("I35").Value = 0.128588910498664
